$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos list refresh diff.
# Some Price (D) values are plain decimal numbers; force them to stay text
# (matching the original inlineStr/text cell type) by briefly switching the
# cell to a text number format, assigning the value, then restoring the
# default "Normal" style so no stray formatting is left behind.

$ws.Range('D2').Value = '42.702.54'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '2.310.19'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '95.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.503'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.29%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0782'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('E14').Value = '  -1.90%  '
$ws.Range('D15').Value = '2.671.62'
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D16').Value = '2.333.98'
$ws.Range('E16').Value = '  +1.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.787'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').Value = '42.654.08'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.12'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.62%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.27'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '235.07'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.25%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.27'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('E28').Value = '  +14.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.53'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('E32').Value = '  +0.01%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.58'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  +1.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0696'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.100'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +15.82%  '
$ws.Range('D43').Value = '1.924.70'
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0278'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('E46').Value = '  -2.13%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('D49').Value = '2.539.62'
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.05'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.69%  '
